$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (column C) date serial from 46075 to 46076 for rows 2-24
$ws.Range("C2:C24").Value = 46076
